# Worked on feature extraction, changed name "previous_tuple" to "tuple_minus1",
# added some global lists for feature extraction.
#
# This adds the "Deviation_from_Biber" notes for a handful of features on the
# "Tabelle2" sheet, and leaves the grid scrolled/selected near the last-edited
# cell (C42), matching the author's on-save view state as closely as possible.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle2")
$ws.Activate()

# Written in the same order the new shared strings appear in the saved file
# (166=none, 167="So far, ...", 168="Biber says ...", 169="The tag we are
# using ...", 170="Here, Biber also counts ...").

# C65 - "phrasal coordination" (feature 64)
$ws.Range("C65").Value = "none"

# C66 - "non-phrasal coordination" (feature 65)
$ws.Range("C66").Value = "So far, the identification of demonstrative pronouns is probably still too imprecise. We would need to re-use the extraction for feature 10. "

# C52 - "demonstratives" (feature 51)
$ws.Range("C52").Value = 'Biber says that here he "excludes demonstrative pronouns (no. 10) and that as relative, complementizer, or subordinator." Instead, we look at all words tagged as determiners and then count the ones that are "that, this, these, those"'

# C34 - "wh-relatives with pied piping" (feature 33)
$ws.Range("C34").Value = "The tag we are using to find prepositions (IN) does probably not overlap 100% with what Biber understood as prepositions"

# C42 - "predicative adjectives" (feature 41)
$ws.Range("C42").Value = 'Here, Biber also counts "any ADJ not identified as predicative - no. 41)", so far we do not do that'

# C62 - "stranded prepositions" (feature 61) - reuses the same shared string as C34
$ws.Range("C62").Value = "The tag we are using to find prepositions (IN) does probably not overlap 100% with what Biber understood as prepositions"

# Leave the view scrolled to row 34 and the active cell on C42, as in the saved file.
$excel.ActiveWindow.ScrollRow = 34
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C42").Select()
